$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 (B:G) with new values ---
$data = @{
    2 = @(0.5418174940949013, 1.620060478675336,  9.761454062984352,  3.124332578805328, 3.110994274217187, 46)
    3 = @(0.5681188649824511, 1.578281505780222,  9.03108104079497,   3.005175708805555, 2.984332153941157, 45)
    4 = @(0.6391235188670518, 1.488536360484932,  9.31892943597502,   3.052692162006353, 3.01954803640462,  44)
    5 = @(0.5562796024915397, 1.588854171992465,  9.477300844941421,  3.078522510059237, 3.063680153555256, 43)
    6 = @(0.7091244686859514, 1.549719313849029,  9.643609441603539,  3.105416146284349, 3.060015865731248, 42)
    7 = @(0.6223279867891005, 1.593628882894858,  9.828395328189679,  3.135027165462794, 3.110808702761217, 41)
    8 = @(0.7758631181458544, 1.561459785857586,  10.03195630587154,  3.167326365544217, 3.10994974446534,  40)
    9 = @(0.5985805415870319, 1.595506958100287,  10.14980463674964,  3.185875803723309, 3.170043719220926, 39)
    10 = @(0.8159597391921622, 1.587015464082368, 10.56404781431108,  3.25023811655563,  3.188381719900689, 38)
    11 = @(0.6067207657392607, 1.576071087915743, 10.50536641133343,  3.241198298674956, 3.227823845440115, 37)
}

# --- Create the new rows (6-11) for column A, copying the format/style used by A2:A5 ---
$ws.Range("A2:A5").Copy()
$ws.Range("A6:A11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$labels = @("Q4", "Q5", "Q6", "Q7", "Q8", "Q9")
for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item(6 + $i, 1).Value = $labels[$i]
}

# --- Write the numeric values for columns B..G on rows 2-11 ---
foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($col = 0; $col -lt $vals.Length; $col++) {
        $ws.Cells.Item($row, 2 + $col).Value = $vals[$col]
    }
}
